$wb = $excel.ActiveWorkbook

$wsKayitlar = $wb.Worksheets.Item("Kayitlar")
$wsKayitlar.Rows.Item(551).Delete()

$wsMerkez = $wb.Worksheets.Item("Merkez İlçe")
$wsMerkez.Rows.Item(7).Delete()
